$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.316.64'
$ws.Range('E2').Value = '  +1.84%  '
$ws.Range('D3').Value = '2.692.07'
$ws.Range('E3').Value = '  +1.74%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '611.57'
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.49'
$ws.Range('E6').Value = '  +2.78%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.592'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +8.83%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.00'
$ws.Range('E10').Value = '  +3.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.406'
$ws.Range('E11').Value = '  +2.02%  '
$ws.Range('E12').Value = '  +1.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000210'
$ws.Range('E13').Value = '  +22.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '30.22'
$ws.Range('E14').Value = '  +4.14%  '
$ws.Range('D15').Value = '3.176.45'
$ws.Range('E15').Value = '  +1.85%  '
$ws.Range('D16').Value = '66.128.50'
$ws.Range('E16').Value = '  +1.87%  '
$ws.Range('D17').Value = '2.692.23'
$ws.Range('E17').Value = '  +1.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.73'
$ws.Range('E18').Value = '  +1.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.92'
$ws.Range('E19').Value = '  +2.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '363.92'
$ws.Range('E20').Value = '  +3.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.55'
$ws.Range('E21').Value = '  +4.32%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.16'
$ws.Range('E23').Value = '  +3.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.76'
$ws.Range('E24').Value = '  +3.09%  '
$ws.Range('B25').Value = 'SuiNetwork'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.68'
$ws.Range('E25').Value = '  -2.51%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000107'
$ws.Range('E26').Value = '  +16.37%  '
$ws.Range('E27').Value = '  +6.15%  '
$ws.Range('E28').Value = '  +0.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.23'
$ws.Range('E29').Value = '  -1.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.22'
$ws.Range('E30').Value = '  +7.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '541.01'
$ws.Range('E31').Value = '  +0.50%  '
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.81'
$ws.Range('E33').Value = '  -0.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.65'
$ws.Range('E34').Value = '  +4.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.46'
$ws.Range('E35').Value = '  -5.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.435'
$ws.Range('E36').Value = '  +1.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.83'
$ws.Range('E37').Value = '  +3.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '163.42'
$ws.Range('E38').Value = '  -1.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.01'
$ws.Range('E39').Value = '  -1.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  +0.34%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '171.63'
$ws.Range('E41').Value = '  +2.22%  '
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('E43').Value = '  +2.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.23'
$ws.Range('E44').Value = '  +2.93%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.35'
$ws.Range('E45').Value = '  +4.27%  '
$ws.Range('E46').Value = '  +2.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.28'
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.665'
$ws.Range('E48').Value = '  +2.89%  '
$ws.Range('E49').Value = '  +6.52%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.45'
$ws.Range('E50').Value = '  +4.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0990'
$ws.Range('E51').Value = '  +0.55%  '
